# Refresh cryptos.xlsx price/volume columns (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.464.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.77%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.604.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.87%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "202.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "594.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.27%  "

$ws.Range("E7").Value = "  +0.74%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +5.93%  "

$ws.Range("E10").Value = "  -0.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.53"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.16%  "

$ws.Range("E12").Value = "  -0.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "698.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +17.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.172.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "70.524.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.75%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.611.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.33%  "

$ws.Range("E20").Value = "  +0.47%  "

$ws.Range("E21").Value = "  +1.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "110.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.69%  "

$ws.Range("E28").Value = "  -0.53%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.91%  "

$ws.Range("E33").Value = "  -0.10%  "

$ws.Range("E34").Value = "  -0.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "63.73"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0851"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.842.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("E38").Value = "  -0.13%  "

$ws.Range("E39").Value = "  -5.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "508.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.16%  "

$ws.Range("E43").Value = "  -1.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.137"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.39%  "

$ws.Range("E45").Value = "  +5.32%  "

$ws.Range("E46").Value = "  +8.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.141"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.91%  "

$ws.Range("E50").Value = "  -0.29%  "

$ws.Range("E51").Value = "  +23.45%  "

